# Update automàtic: dades i banners [2026-02-28 19:19]
# Refreshes DATA_EXTRACCIO timestamps and the latest observation values
# scraped from meteo.cat for the affected station rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-28 19:18:13'
$ws.Range("E3").Value = '2026-02-28 19:18:16'
$ws.Range("I3").Value = '0.3 mm'
$ws.Range("N3").Value = '-2.9 °C 18:59 TU'
$ws.Range("E4").Value = '2026-02-28 19:18:18'
$ws.Range("H4").Value = '''82%'
$ws.Range("O4").Value = '11.2 °C'
$ws.Range("E5").Value = '2026-02-28 19:18:20'
$ws.Range("N5").Value = '-2.7 °C 18:55 TU'
$ws.Range("E6").Value = '2026-02-28 19:18:23'
$ws.Range("E7").Value = '2026-02-28 19:18:25'
$ws.Range("E8").Value = '2026-02-28 19:18:28'
$ws.Range("E9").Value = '2026-02-28 19:18:30'
$ws.Range("H9").Value = '''80%'
$ws.Range("E10").Value = '2026-02-28 19:18:31'
$ws.Range("E11").Value = '2026-02-28 19:18:32'
$ws.Range("O11").Value = '7.2 °C'
$ws.Range("E12").Value = '2026-02-28 19:18:33'
$ws.Range("E13").Value = '2026-02-28 19:18:34'
$ws.Range("H13").Value = '''68%'
$ws.Range("J13").Value = '1024.1 hPa'
$ws.Range("E14").Value = '2026-02-28 19:18:35'
$ws.Range("E15").Value = '2026-02-28 19:18:36'
$ws.Range("O15").Value = '11.2 °C'
$ws.Range("E16").Value = '2026-02-28 19:18:38'
$ws.Range("N16").Value = '-2.7 °C 18:59 TU'
$ws.Range("E17").Value = '2026-02-28 19:18:39'
$ws.Range("H17").Value = '''81%'
$ws.Range("O17").Value = '3.0 °C'
$ws.Range("E18").Value = '2026-02-28 19:18:40'
$ws.Range("E19").Value = '2026-02-28 19:18:41'
$ws.Range("H19").Value = '''77%'
$ws.Range("O19").Value = '7.9 °C'
$ws.Range("E20").Value = '2026-02-28 19:18:43'
$ws.Range("H20").Value = '''61%'
$ws.Range("N20").Value = '-1.9 °C 18:55 TU'
$ws.Range("E21").Value = '2026-02-28 19:18:46'
$ws.Range("E22").Value = '2026-02-28 19:18:48'
$ws.Range("H22").Value = '''67%'
$ws.Range("N22").Value = '-2.9 °C 18:50 TU'
$ws.Range("E23").Value = '2026-02-28 19:18:51'
$ws.Range("H23").Value = '''69%'
$ws.Range("I23").Value = '0.1 mm'
$ws.Range("N23").Value = '-2.1 °C 18:47 TU'
$ws.Range("E24").Value = '2026-02-28 19:18:53'
$ws.Range("J24").Value = '1025.1 hPa'
$ws.Range("E25").Value = '2026-02-28 19:18:55'
$ws.Range("H25").Value = '''60%'
$ws.Range("N25").Value = '-0.9 °C 18:37 TU'
$ws.Range("E26").Value = '2026-02-28 19:18:58'
$ws.Range("J26").Value = '1024.3 hPa'
$ws.Range("O26").Value = '4.9 °C'
$ws.Range("E27").Value = '2026-02-28 19:19:00'
$ws.Range("H27").Value = '''53%'
$ws.Range("N27").Value = '-0.5 °C 18:56 TU'
$ws.Range("O27").Value = '2.0 °C'
$ws.Range("E28").Value = '2026-02-28 19:19:03'
$ws.Range("E29").Value = '2026-02-28 19:19:05'
$ws.Range("E30").Value = '2026-02-28 19:19:07'
$ws.Range("E31").Value = '2026-02-28 19:19:10'
$ws.Range("K31").Value = '11.9 MJ/m2'
$ws.Range("E32").Value = '2026-02-28 19:19:12'
$ws.Range("E33").Value = '2026-02-28 19:19:15'
$ws.Range("J33").Value = '1022.9 hPa'
$ws.Range("O33").Value = '7.3 °C'
$ws.Range("E34").Value = '2026-02-28 19:19:17'
$ws.Range("H34").Value = '''66%'
$ws.Range("E35").Value = '2026-02-28 19:19:20'
$ws.Range("E36").Value = '2026-02-28 19:19:22'
$ws.Range("E37").Value = '2026-02-28 19:19:24'
$ws.Range("E38").Value = '2026-02-28 19:19:27'
$ws.Range("E39").Value = '2026-02-28 19:19:29'
$ws.Range("H39").Value = '''61%'
$ws.Range("N39").Value = '-2.0 °C 18:56 TU'
$ws.Range("O39").Value = '-0.5 °C'
$ws.Range("E40").Value = '2026-02-28 19:19:31'
$ws.Range("H40").Value = '''75%'
$ws.Range("J40").Value = '1024.5 hPa'
$ws.Range("E41").Value = '2026-02-28 19:19:33'
$ws.Range("J41").Value = '1024.5 hPa'
$ws.Range("E42").Value = '2026-02-28 19:19:36'
$ws.Range("E43").Value = '2026-02-28 19:19:38'
$ws.Range("H43").Value = '''80%'
$ws.Range("O43").Value = '7.6 °C'
$ws.Range("E44").Value = '2026-02-28 19:19:40'
$ws.Range("I44").Value = '2.7 mm'
$ws.Range("O44").Value = '-1.0 °C'
$ws.Range("E45").Value = '2026-02-28 19:19:43'
$ws.Range("I45").Value = '1.2 mm'
$ws.Range("N45").Value = '4.4 °C 18:52 TU'
$ws.Range("E46").Value = '2026-02-28 19:19:45'
$ws.Range("N46").Value = '9.9 °C 18:59 TU'
$ws.Range("O46").Value = '11.6 °C'
